$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '41.233.01'
$ws.Range("E2").Value = '  -1.54%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.176.74'
$ws.Range("E3").Value = '  -2.30%  '

$ws.Range("E4").Value = '  -0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '248.13'
$ws.Range("E5").Value = '  -1.06%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.611'
$ws.Range("E6").Value = '  -2.90%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '65.72'
$ws.Range("E7").Value = '  -8.19%  '

$ws.Range("E8").Value = '  +0.05%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.564'
$ws.Range("E9").Value = '  -4.29%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '59.66'
$ws.Range("E10").Value = '  +2.54%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0926'
$ws.Range("E11").Value = '  -4.55%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '35.64'
$ws.Range("E12").Value = '  -13.42%  '

$ws.Range("E13").Value = '  -1.40%  '

$ws.Range("E14").Value = '  -4.20%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.502.73'
$ws.Range("E15").Value = '  -2.26%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.32'
$ws.Range("E16").Value = '  -4.35%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.848'
$ws.Range("E17").Value = '  -2.06%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.177.01'
$ws.Range("E18").Value = '  -2.38%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '41.110.17'
$ws.Range("E19").Value = '  -1.54%  '

$ws.Range("E20").Value = '  -3.56%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.08'
$ws.Range("E21").Value = '  -2.26%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '71.38'
$ws.Range("E22").Value = '  -2.20%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '229.32'
$ws.Range("E23").Value = '  -2.59%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.05'
$ws.Range("E24").Value = '  -4.47%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.88'
$ws.Range("E25").Value = '  -6.60%  '

$ws.Range("E26").Value = '  +0.14%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.25'
$ws.Range("E27").Value = '  +5.31%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.41'
$ws.Range("E28").Value = '  -4.68%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.72'
$ws.Range("E29").Value = '  -5.69%  '

$ws.Range("E30").Value = '  -3.47%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '167.97'
$ws.Range("E31").Value = '  -1.84%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.18'
$ws.Range("E32").Value = '  -2.80%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.122'
$ws.Range("E33").Value = '  -0.59%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.70'
$ws.Range("E34").Value = '  +0.98%  '

$ws.Range("E35").Value = '  +2.31%  '

$ws.Range("E36").Value = '  -3.78%  '

$ws.Range("E37").Value = '  -3.85%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.94'
$ws.Range("E38").Value = '  +0.47%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '24.03'
$ws.Range("E39").Value = '  -8.25%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0303'
$ws.Range("E40").Value = '  +0.11%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.19'
$ws.Range("E41").Value = '  -4.99%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.45'
$ws.Range("E42").Value = '  -8.73%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.91'
$ws.Range("E43").Value = '  +3.30%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '60.51'
$ws.Range("E44").Value = '  -10.81%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '11.08'
$ws.Range("E45").Value = '  -6.16%  '

$ws.Range("E46").Value = '  -8.58%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.47'
$ws.Range("E47").Value = '  -4.27%  '

$ws.Range("E48").Value = '  -0.09%  '

$ws.Range("E49").Value = '  -3.29%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.15'
$ws.Range("E50").Value = '  -1.13%  '

$ws.Range("E51").Value = '  -4.56%  '
